$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.620.81'
$ws.Range("D3").Value = '2.270.08'
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '118.57'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +7.94%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '267.01'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.35%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.642'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.13%  '
$ws.Range("E8").Value = '  +0.20%  '
$ws.Range("E9").Value = '  +3.31%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '47.33'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.63%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0943'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.88%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '9.47'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +8.93%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.69'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.33%  '
$ws.Range("E15").Value = '  +6.32%  '
$ws.Range("D16").Value = '2.613.10'
$ws.Range("E16").Value = '  +0.29%  '
$ws.Range("D17").Value = '2.268.84'
$ws.Range("E17").Value = '  +0.12%  '
$ws.Range("D18").Value = '43.639.69'
$ws.Range("E18").Value = '  +1.56%  '
$ws.Range("E19").Value = '  +2.42%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.91'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.88%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.35'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.00%  '
$ws.Range("E22").Value = '  -2.72%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '234.21'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.66%  '
$ws.Range("E24").Value = '  +3.24%  '
$ws.Range("E25").Value = '  -0.42%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.25'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +9.12%  '
$ws.Range("E27").Value = '  +1.92%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '42.04'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.79%  '
$ws.Range("E29").Value = '  +2.23%  '
$ws.Range("E30").Value = '  +0.32%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '174.51'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.90%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.46'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.27%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0914'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.84%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.73'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.76%  '
$ws.Range("E35").Value = '  +3.00%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.26'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +12.83%  '
$ws.Range("E37").Value = '  +9.76%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.59'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.78%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.109'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.45%  '
$ws.Range("E40").Value = '  -1.16%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '13.87'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.57%  '
$ws.Range("E42").Value = '  +2.67%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '72.49'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.80%  '
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.38'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.65%  '
$ws.Range("E46").Value = '  -4.82%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '75.61'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +43.48%  '
$ws.Range("B48").Value = 'TheSandbox'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.665'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +20.09%  '
$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.57'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.18%  '
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '103.01'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.95%  '
$ws.Range("B51").Value = 'TrustWalletToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.26'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.28%  '
